$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: cell reference -> new value, as scraped from the updated
# coinranking.com snapshot (GitHub Actions run on 2024-03-02).
$updates = @(
    @{ Cell = "D2"; Value = '62.142.84' }
    @{ Cell = "E2"; Value = '  +0.86%  ' }
    @{ Cell = "D3"; Value = '3.443.56' }
    @{ Cell = "E3"; Value = '  +1.62%  ' }
    @{ Cell = "E4"; Value = '  +0.04%  ' }
    @{ Cell = "D5"; Value = '409.68' }
    @{ Cell = "E5"; Value = '  +1.08%  ' }
    @{ Cell = "D6"; Value = '129.21' }
    @{ Cell = "E6"; Value = '  -3.26%  ' }
    @{ Cell = "D7"; Value = '0.625' }
    @{ Cell = "E7"; Value = '  +5.86%  ' }
    @{ Cell = "D9"; Value = '0.741' }
    @{ Cell = "E9"; Value = '  +10.86%  ' }
    @{ Cell = "E10"; Value = '  +16.53%  ' }
    @{ Cell = "D11"; Value = '42.95' }
    @{ Cell = "E11"; Value = '  +1.15%  ' }
    @{ Cell = "E12"; Value = '  -0.36%  ' }
    @{ Cell = "D13"; Value = '3.978.23' }
    @{ Cell = "E13"; Value = '  +1.51%  ' }
    @{ Cell = "D14"; Value = '21.34' }
    @{ Cell = "E14"; Value = '  +8.00%  ' }
    @{ Cell = "D15"; Value = '8.96' }
    @{ Cell = "E15"; Value = '  +6.61%  ' }
    @{ Cell = "D16"; Value = '0.0000207' }
    @{ Cell = "E16"; Value = '  +61.56%  ' }
    @{ Cell = "D17"; Value = '3.481.51' }
    @{ Cell = "E17"; Value = '  +2.79%  ' }
    @{ Cell = "D18"; Value = '12.44' }
    @{ Cell = "E18"; Value = '  +12.56%  ' }
    @{ Cell = "E19"; Value = '  +5.45%  ' }
    @{ Cell = "D20"; Value = '62.164.08' }
    @{ Cell = "E20"; Value = '  +1.01%  ' }
    @{ Cell = "D21"; Value = '407.77' }
    @{ Cell = "E21"; Value = '  +28.88%  ' }
    @{ Cell = "D22"; Value = '90.11' }
    @{ Cell = "E22"; Value = '  +5.37%  ' }
    @{ Cell = "D23"; Value = '3.20' }
    @{ Cell = "E23"; Value = '  +0.10%  ' }
    @{ Cell = "D24"; Value = '13.47' }
    @{ Cell = "E24"; Value = '  +5.35%  ' }
    @{ Cell = "E25"; Value = '  +2.85%  ' }
    @{ Cell = "D26"; Value = '33.37' }
    @{ Cell = "E26"; Value = '  +13.07%  ' }
    @{ Cell = "D27"; Value = '8.74' }
    @{ Cell = "E27"; Value = '  +5.20%  ' }
    @{ Cell = "E28"; Value = '  +0.39%  ' }
    @{ Cell = "B29"; Value = 'Toncoin' }
    @{ Cell = "C29"; Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton' }
    @{ Cell = "D29"; Value = '2.79' }
    @{ Cell = "E29"; Value = '  +5.69%  ' }
    @{ Cell = "B30"; Value = 'RenderToken' }
    @{ Cell = "C30"; Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr' }
    @{ Cell = "D30"; Value = '7.61' }
    @{ Cell = "E30"; Value = '  +0.42%  ' }
    @{ Cell = "E31"; Value = '  +2.23%  ' }
    @{ Cell = "B32"; Value = 'Cosmos' }
    @{ Cell = "C32"; Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom' }
    @{ Cell = "D32"; Value = '11.93' }
    @{ Cell = "E32"; Value = '  +5.14%  ' }
    @{ Cell = "B33"; Value = 'Kaspa' }
    @{ Cell = "C33"; Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas' }
    @{ Cell = "D33"; Value = '0.172' }
    @{ Cell = "E33"; Value = '  +0.45%  ' }
    @{ Cell = "D34"; Value = '43.74' }
    @{ Cell = "E34"; Value = '  +6.52%  ' }
    @{ Cell = "E35"; Value = '  -0.04%  ' }
    @{ Cell = "E36"; Value = '  +5.36%  ' }
    @{ Cell = "D37"; Value = '54.35' }
    @{ Cell = "E37"; Value = '  +5.05%  ' }
    @{ Cell = "D38"; Value = '0.999' }
    @{ Cell = "E38"; Value = '  +0.10%  ' }
    @{ Cell = "D39"; Value = '3.40' }
    @{ Cell = "E39"; Value = '  -0.46%  ' }
    @{ Cell = "E40"; Value = '  -0.17%  ' }
    @{ Cell = "E41"; Value = '  +6.68%  ' }
    @{ Cell = "D42"; Value = '0.314' }
    @{ Cell = "E42"; Value = '  +6.42%  ' }
    @{ Cell = "D43"; Value = '141.38' }
    @{ Cell = "E43"; Value = '  +1.30%  ' }
    @{ Cell = "D44"; Value = '1.99' }
    @{ Cell = "E44"; Value = '  +0.87%  ' }
    @{ Cell = "D45"; Value = '4.05' }
    @{ Cell = "E45"; Value = '  +1.59%  ' }
    @{ Cell = "E46"; Value = '  +7.57%  ' }
    @{ Cell = "D47"; Value = '16.79' }
    @{ Cell = "E47"; Value = '  +1.15%  ' }
    @{ Cell = "D48"; Value = '22.06' }
    @{ Cell = "E48"; Value = '  +4.08%  ' }
    @{ Cell = "D49"; Value = '2.124.39' }
    @{ Cell = "E49"; Value = '  +0.33%  ' }
    @{ Cell = "D50"; Value = '0.129' }
    @{ Cell = "E50"; Value = '  +14.97%  ' }
    @{ Cell = "B51"; Value = 'BEAM' }
    @{ Cell = "C51"; Value = 'https://coinranking.com/coin/cYYMfXF4u+beam-beam' }
    @{ Cell = "D51"; Value = '0.0376' }
    @{ Cell = "E51"; Value = '  +7.46%  ' }
)

foreach ($u in $updates) {
    $range = $ws.Range($u.Cell)
    if ($u.Cell.StartsWith("D")) {
        # Column D holds price text that can look like a number (e.g. "409.68")
        # or like a dotted-thousands price (e.g. "62.142.84"). Force text format
        # so Excel does not silently convert it to a numeric value.
        $range.NumberFormat = "@"
    }
    $range.Value = $u.Value
}
